# Commit: "Add challenges and needs"
#
# Inserts two new slides right after slide 2 ("Objectives") and before the
# former slide 3 ("Start at the end!  What is your end goal?"):
#   - "Challenges of Training for XebiaLabs Tools"
#   - "Needs of our Users"
#
# Both new slides use the "Title and Content" layout (same layout already
# used by the surrounding content slides).

$p = $ppt.ActivePresentation

# "Title and Content" is the 2nd custom layout on the slide master.
$titleAndContent = $p.SlideMaster.CustomLayouts.Item(2)

# --- New slide 3: "Challenges of Training for XebiaLabs Tools" -----------
$s3 = $p.Slides.AddSlide(3, $titleAndContent)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Challenges of Training for XebiaLabs Tools"

$s3Body = $s3.Shapes.Item(2).TextFrame.TextRange
$s3Body.Text = "`rNew model-based concept`rComplex products`rMany interfaces to other complex products:  middleware, artifact repositories, security, etc.`rWide range of user skill levels`r"

# --- New slide 4: "Needs of our Users" ------------------------------------
$s4 = $p.Slides.AddSlide(4, $titleAndContent)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Needs of our Users"

$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
$s4Body.Text = "Basic:  running a deployment or a release`rAdvanced:  scripts, plugins, interfacing with API`rAdmin 1:  high-availability, databases, security`rAdmin 2: structuring folders and permissions`r"
